$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 updates
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 updates: B2, D2, E2 cleared; C2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 24.716124480553788
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 updates
$ws.Range("B3").Value = 22.129156949644027
$ws.Range("C3").Value = 28.037104532165927
$ws.Range("D3").Value = 28.046525003550357
$ws.Range("E3").Value = 12.731089505263753

# Update selection to match new range
$ws.Range("B1:E3").Select()
